# Generate Report for Archive
# Swap the "fead0616..." and "268b50a0..." records on every sheet (rows 6/7),
# moving fead0616's record to row 6 and 268b50a0's record to row 7, and
# changing the (new) row 6 Status from "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" (columns A,B,C ; hyperlinks on column A)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A6").Value = "fead0616-771c-481d-b56a-5918df0efd59.md"
$ws1.Range("B6").Value = "In Translation"
$ws1.Range("C6").Value = "In Translation"

$ws1.Range("A7").Value = "268b50a0-f412-4f69-99e2-079bdfdf1585.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md")
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" (columns A,B,C,D ; hyperlinks on columns A and C)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A6").Value = "fead0616-771c-481d-b56a-5918df0efd59.md"
$ws2.Range("B6").Value = "In Translation"
$ws2.Range("C6").Value = "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.zh-cn.xlf"
$ws2.Range("D6").Value = "2016-01-25 06:14:58"

$ws2.Range("A7").Value = "268b50a0-f412-4f69-99e2-079bdfdf1585.md"
$ws2.Range("B7").Value = "Ready for handoff"
$ws2.Range("C7").Value = "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.zh-cn.xlf"
$ws2.Range("D7").Value = "2016-01-25 06:10:10"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/288a598b1ee77e39219960e17f80572c35dcfff9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.zh-cn.xlf", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.zh-cn.xlf", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.zh-cn.xlf", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9572f2523077644a3dbe8565fbbe0ddadaeb7d8d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b474860421311042c2e35d5037bf8c00eba3310c/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md")
$ws2.Hyperlinks.Add($ws2.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61a67c1bb3d0bfc0b09208ed5006a571e79f7884/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md")
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08011f83eaae42323656547480ff840f8295f6ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.zh-cn.xlf", "", "", "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a81e18e82e375f9e0868ffb2cff3831b7ecd46ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.zh-cn.xlf", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" (columns A,B,C,D ; hyperlinks on columns A and C)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A6").Value = "fead0616-771c-481d-b56a-5918df0efd59.md"
$ws3.Range("B6").Value = "In Translation"
$ws3.Range("C6").Value = "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.de-de.xlf"
$ws3.Range("D6").Value = "2016-01-25 06:15:12"

$ws3.Range("A7").Value = "268b50a0-f412-4f69-99e2-079bdfdf1585.md"
$ws3.Range("B7").Value = "Ready for handoff"
$ws3.Range("C7").Value = "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.de-de.xlf"
$ws3.Range("D7").Value = "2016-01-25 06:10:21"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20f8d77436701dde648d700b38e617bb690aeba4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.de-de.xlf", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd3ba71bf9b68dfda636a4e40592cec9e35685db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.de-de.xlf", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd3ba71bf9b68dfda636a4e40592cec9e35685db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.de-de.xlf", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92e52a09fc3009b2982b8fcef146c17277ba9ecf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0f827d899564115b2cc4e0f074de84f7815a3845/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md")
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3c70eca77d8230e5c02b9d75fb9a273073453ba7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md")
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a22482c96d73fcaa46b3f0da45e0cc393a55b277/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.de-de.xlf", "", "", "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md")
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63a3f0b36de17ff4353a23d9718d6ac5a1ba2980/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.de-de.xlf", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/435753bad1342968cb141d1e05c3a3d59fbda92d/.localization-config", "", "", ".localization-config")
